# Applies the crypto-tracker refresh described in the commit:
#   "Updated cryptos list on Tue Feb  6 20:23:49 UTC 2024 with GitHub Actions"
#
# Most edits are simple Price (column D) / Volume 1h (column E) refreshes.
# Rows 12 and 13 additionally swap rank position (TRON now above Chainlink).
#
# Column D sometimes holds digit-only strings (e.g. "97.58"). The sheet stores
# these as text, so a leading "'" forces Excel to keep them as text instead of
# silently re-typing the cell as a Number (which would also drop the exact
# decimal formatting, e.g. "0.0730" -> 0.073).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.188.71'
$ws.Range('E2').Value = '  +1.35%  '

# Row 3
$ws.Range('D3').Value = '2.389.49'
$ws.Range('E3').Value = '  +4.21%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').Value = '''303.37'
$ws.Range('E5').Value = '  +0.70%  '

# Row 6
$ws.Range('D6').Value = '''97.58'
$ws.Range('E6').Value = '  +2.19%  '

# Row 7
$ws.Range('D7').Value = '''0.509'
$ws.Range('E7').Value = '  +0.59%  '

# Row 8
$ws.Range('E8').Value = '  -0.12%  '

# Row 9
$ws.Range('E9').Value = '  +2.57%  '

# Row 10
$ws.Range('D10').Value = '''34.33'
$ws.Range('E10').Value = '  +0.06%  '

# Row 11
$ws.Range('D11').Value = '''0.0791'
$ws.Range('E11').Value = '  +0.69%  '

# Row 12
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '''0.122'
$ws.Range('E12').Value = '  +3.00%  '

# Row 13
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = '''18.53'
$ws.Range('E13').Value = '  -2.27%  '

# Row 14
$ws.Range('E14').Value = '  +0.72%  '

# Row 15
$ws.Range('D15').Value = '2.756.36'
$ws.Range('E15').Value = '  +3.69%  '

# Row 16
$ws.Range('D16').Value = '2.378.69'
$ws.Range('E16').Value = '  +3.81%  '

# Row 17
$ws.Range('D17').Value = '''0.816'
$ws.Range('E17').Value = '  +4.81%  '

# Row 18
$ws.Range('D18').Value = '43.189.92'
$ws.Range('E18').Value = '  +1.51%  '

# Row 19
$ws.Range('D19').Value = '''12.29'
$ws.Range('E19').Value = '  +1.24%  '

# Row 20
$ws.Range('E20').Value = '  +6.48%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0891'
$ws.Range('E21').Value = '  +0.41%  '

# Row 22
$ws.Range('D22').Value = '''68.49'
$ws.Range('E22').Value = '  +1.21%  '

# Row 23
$ws.Range('D23').Value = '''236.55'
$ws.Range('E23').Value = '  +0.61%  '

# Row 24
$ws.Range('D24').Value = '''2.23'
$ws.Range('E24').Value = '  -1.50%  '

# Row 26
$ws.Range('E26').Value = '  -0.03%  '

# Row 27
$ws.Range('D27').Value = '''24.89'

# Row 28
$ws.Range('E28').Value = '  +0.18%  '

# Row 29
$ws.Range('D29').Value = '''9.14'
$ws.Range('E29').Value = '  +1.38%  '

# Row 30
$ws.Range('D30').Value = '''31.63'
$ws.Range('E30').Value = '  -0.37%  '

# Row 31
$ws.Range('E31').Value = '  +3.26%  '

# Row 32
$ws.Range('E32').Value = '  -0.04%  '

# Row 33
$ws.Range('D33').Value = '''0.0730'
$ws.Range('E33').Value = '  +4.68%  '

# Row 34
$ws.Range('D34').Value = '''17.28'
$ws.Range('E34').Value = '  -1.07%  '

# Row 35
$ws.Range('E35').Value = '  +7.62%  '

# Row 36
$ws.Range('D36').Value = '''4.39'
$ws.Range('E36').Value = '  -0.27%  '

# Row 38
$ws.Range('E38').Value = '  +2.07%  '

# Row 39
$ws.Range('D39').Value = '''2.82'
$ws.Range('E39').Value = '  +5.64%  '

# Row 40
$ws.Range('D40').Value = '''22.74'
$ws.Range('E40').Value = '  +12.92%  '

# Row 41
$ws.Range('E41').Value = '  +0.68%  '

# Row 42
$ws.Range('D42').Value = '''106.06'
$ws.Range('E42').Value = '  -35.98%  '

# Row 43
$ws.Range('D43').Value = '1.950.18'
$ws.Range('E43').Value = '  -0.63%  '

# Row 44
$ws.Range('E44').Value = '  +1.00%  '

# Row 45
$ws.Range('E45').Value = '  +2.40%  '

# Row 46
$ws.Range('D46').Value = '''9.35'
$ws.Range('E46').Value = '  -10.27%  '

# Row 47
$ws.Range('D47').Value = '''2.76'
$ws.Range('E47').Value = '  +0.50%  '

# Row 48
$ws.Range('D48').Value = '2.617.28'
$ws.Range('E48').Value = '  +3.78%  '

# Row 49
$ws.Range('D49').Value = '''53.11'
$ws.Range('E49').Value = '  +0.25%  '

# Row 50
$ws.Range('D50').Value = '''72.40'
$ws.Range('E50').Value = '  +1.77%  '

# Row 51
$ws.Range('E51').Value = '  +2.03%  '
